$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 40.75339133333333
$ws.Range("H2").Value = 122.260174
$ws.Range("I2").Value = 0.02126536631186857
$ws.Range("J2").Value = 0.02126536631186857
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.776574666666666
$ws.Range("N2").Value = 11.329724
$ws.Range("O2").Value = 0.9855052394405499
$ws.Range("P2").Value = 0.9855052394405499
$ws.Range("Q2").Value = 153.9082252902195
$ws.Range("R2").Value = 1385.174027611976
$ws.Range("S2").Value = 0.02095712991896903
$ws.Range("T2").Value = 0.02095712991896903

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 40.75339133333333
$ws.Range("H3").Value = 122.260174
$ws.Range("I3").Value = 0.02126536631186857
$ws.Range("J3").Value = 0.02126536631186857
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.05554566666666667
$ws.Range("N3").Value = 0.166637
$ws.Range("O3").Value = 0.01449476055945007
$ws.Range("P3").Value = 0.01449476055945008
$ws.Range("Q3").Value = 2.263674290537556
$ws.Range("R3").Value = 20.373068614838
$ws.Range("S3").Value = 0.0003082363928995308
$ws.Range("T3").Value = 0.0003082363928995308

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1689.289306666667
$ws.Range("H4").Value = 5067.86792
$ws.Range("I4").Value = 0.8814813868902838
$ws.Range("J4").Value = 0.8814813868902838
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 3.776574666666666
$ws.Range("N4").Value = 11.329724
$ws.Range("O4").Value = 0.9855052394405499
$ws.Range("P4").Value = 0.9855052394405499
$ws.Range("Q4").Value = 6379.727200228231
$ws.Range("R4").Value = 57417.54480205407
$ws.Range("S4").Value = 0.8687045252496971
$ws.Range("T4").Value = 0.8687045252496971

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1689.289306666667
$ws.Range("H5").Value = 5067.86792
$ws.Range("I5").Value = 0.8814813868902838
$ws.Range("J5").Value = 0.8814813868902838
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.05554566666666667
$ws.Range("N5").Value = 0.166637
$ws.Range("O5").Value = 0.01449476055945007
$ws.Range("P5").Value = 0.01449476055945008
$ws.Range("Q5").Value = 93.83270073167111
$ws.Range("R5").Value = 844.49430658504
$ws.Range("S5").Value = 0.01277686164058664
$ws.Range("T5").Value = 0.01277686164058664

# Row 6
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Fn1"
$ws.Range("C6").Value = "Tnfrsf11b"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 100.9654023333333
$ws.Range("H6").Value = 302.896207
$ws.Range("I6").Value = 0.05268435816499466
$ws.Range("J6").Value = 0.05268435816499466
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.776574666666666
$ws.Range("N6").Value = 11.329724
$ws.Range("O6").Value = 0.9855052394405499
$ws.Range("P6").Value = 0.9855052394405499
$ws.Range("Q6").Value = 381.3033806618742
$ws.Range("R6").Value = 3431.730425956868
$ws.Range("S6").Value = 0.05192071100816475
$ws.Range("T6").Value = 0.05192071100816475

# Row 7
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Fn1"
$ws.Range("C7").Value = "Tnfrsf11b"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 100.9654023333333
$ws.Range("H7").Value = 302.896207
$ws.Range("I7").Value = 0.05268435816499466
$ws.Range("J7").Value = 0.05268435816499466
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.05554566666666667
$ws.Range("N7").Value = 0.166637
$ws.Range("O7").Value = 0.01449476055945007
$ws.Range("P7").Value = 0.01449476055945008
$ws.Range("Q7").Value = 5.608190582873222
$ws.Range("R7").Value = 50.473715245859
$ws.Range("S7").Value = 0.0007636471568299061
$ws.Range("T7").Value = 0.0007636471568299062

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Fn1"
$ws.Range("C8").Value = "Tnfrsf11b"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 85.41274733333334
$ws.Range("H8").Value = 256.238242
$ws.Range("I8").Value = 0.04456888863285297
$ws.Range("J8").Value = 0.04456888863285297
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.776574666666666
$ws.Range("N8").Value = 11.329724
$ws.Range("O8").Value = 0.9855052394405499
$ws.Range("P8").Value = 0.9855052394405499
$ws.Range("Q8").Value = 322.5676177894676
$ws.Range("R8").Value = 2903.108560105208
$ws.Range("S8").Value = 0.04392287326371896
$ws.Range("T8").Value = 0.04392287326371896

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Fn1"
$ws.Range("C9").Value = "Tnfrsf11b"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 85.41274733333334
$ws.Range("H9").Value = 256.238242
$ws.Range("I9").Value = 0.04456888863285297
$ws.Range("J9").Value = 0.04456888863285297
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.05554566666666667
$ws.Range("N9").Value = 0.166637
$ws.Range("O9").Value = 0.01449476055945007
$ws.Range("P9").Value = 0.01449476055945008
$ws.Range("Q9").Value = 4.744307992461556
$ws.Range("R9").Value = 42.698771932154
$ws.Range("S9").Value = 0.0006460153691339998
$ws.Range("T9").Value = 0.000646015369134

